$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header cells: "_old" columns become "_FV2404", "_new" columns become "_FV2410".
$headerRenames = @(
    @{ Cell = "A1"; Value = "Segmentname_FV2404" },
    @{ Cell = "B1"; Value = "Segmentgruppe_FV2404" },
    @{ Cell = "C1"; Value = "Segment_FV2404" },
    @{ Cell = "D1"; Value = "Datenelement_FV2404" },
    @{ Cell = "E1"; Value = "Segment ID_FV2404" },
    @{ Cell = "F1"; Value = "Code_FV2404" },
    @{ Cell = "G1"; Value = "Qualifier_FV2404" },
    @{ Cell = "H1"; Value = "Beschreibung_FV2404" },
    @{ Cell = "I1"; Value = "Bedingungsausdruck_FV2404" },
    @{ Cell = "J1"; Value = "Bedingung_FV2404" },
    @{ Cell = "L1"; Value = "Segmentname_FV2410" },
    @{ Cell = "M1"; Value = "Segmentgruppe_FV2410" },
    @{ Cell = "N1"; Value = "Segment_FV2410" },
    @{ Cell = "O1"; Value = "Datenelement_FV2410" },
    @{ Cell = "P1"; Value = "Segment ID_FV2410" },
    @{ Cell = "Q1"; Value = "Code_FV2410" },
    @{ Cell = "R1"; Value = "Qualifier_FV2410" },
    @{ Cell = "S1"; Value = "Beschreibung_FV2410" },
    @{ Cell = "T1"; Value = "Bedingungsausdruck_FV2410" },
    @{ Cell = "U1"; Value = "Bedingung_FV2410" }
)

foreach ($rename in $headerRenames) {
    $ws.Range($rename.Cell).Value = $rename.Value
}

# 2) Turn the header/data range into a real table (adds xl/tables/table1.xml,
#    the table relationship and <tableParts> on the sheet, plus the AutoFilter).
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U55"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3) Freeze the header row (split beneath row 1, top-left of the scrolling pane is A2).
$ws.Range("A2").Activate()
$excel.ActiveWindow.FreezePanes = $true
